# Weekly fruit/vegetable price update: insert this week's new "Ajo" (Chino,
# Primera) record for Vega Monumental Concepción just after the most recent
# existing entry (row 265), pushing the older history rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 266:277 down to 267:278, opening up a blank row 266.
$ws.Rows.Item(266).Insert()

# Populate the new row 266 with this week's data.
$ws.Range("A266").Value = 11
$ws.Range("B266").Value = "Vega Monumental Concepción"
$ws.Range("C266").Value = "Bíobío"
$ws.Range("D266").Value = 45041
$ws.Range("E266").Value = 8
$ws.Range("F266").Value = 100112003
$ws.Range("G266").Value = "Ajo"
$ws.Range("H266").Value = "Chino"
$ws.Range("I266").Value = "Primera"
$ws.Range("J266").Value = 220
$ws.Range("K266").Value = 14000
$ws.Range("L266").Value = 15000
$ws.Range("M266").Value = 14455
$ws.Range("N266").Value = "$/caja 10 kilos"
$ws.Range("O266").Value = "China"
$ws.Range("P266").Value = 1446
$ws.Range("Q266").Value = 10
$ws.Range("R266").Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Range("D266").NumberFormat = "YYYY-MM-DD HH:MM:SS"
